$wb = $excel.ActiveWorkbook

# ALC row 51
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(51, 8).Value = 0
$ws.Cells.Item(51, 9).Value = 0
$ws.Cells.Item(51, 11).Value = 0
$ws.Cells.Item(51, 13).Value = $null

# ALC row 69
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(69, 8).Value = 7463.8184
$ws.Cells.Item(69, 9).Value = 7004.6665
$ws.Cells.Item(69, 10).Value = 8014.8
$ws.Cells.Item(69, 11).Value = 21013.9995
$ws.Cells.Item(69, 12).Value = 24044.4
$ws.Cells.Item(69, 13).Value = -20139.9995
$ws.Cells.Item(69, 14).Value = -25792.4

# ALC row 72
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(72, 8).Value = 7463.8184
$ws.Cells.Item(72, 9).Value = 7004.6665
$ws.Cells.Item(72, 10).Value = 8014.8
$ws.Cells.Item(72, 11).Value = 63041.9985
$ws.Cells.Item(72, 12).Value = 72133.2
$ws.Cells.Item(72, 13).Value = -58673.9985
$ws.Cells.Item(72, 14).Value = -80869.2

# ALC row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 2207.3
$ws.Cells.Item(98, 9).Value = 1646.75
$ws.Cells.Item(98, 10).Value = 4449.5
$ws.Cells.Item(98, 11).Value = 1646.75
$ws.Cells.Item(98, 12).Value = 4449.5
$ws.Cells.Item(98, 13).Value = -148.75
$ws.Cells.Item(98, 14).Value = -7445.5

# ALC row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(122, 8).Value = 2207.3
$ws.Cells.Item(122, 9).Value = 1646.75
$ws.Cells.Item(122, 10).Value = 4449.5
$ws.Cells.Item(122, 11).Value = 4940.25
$ws.Cells.Item(122, 12).Value = 13348.5
$ws.Cells.Item(122, 13).Value = -2490.25
$ws.Cells.Item(122, 14).Value = -18248.5

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 2372.125
$ws.Cells.Item(137, 9).Value = 2329.8333
$ws.Cells.Item(137, 11).Value = 6989.499899999999
$ws.Cells.Item(137, 13).Value = -4439.499899999999

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5140
$ws.Cells.Item(32, 9).Value = 5925
$ws.Cells.Item(32, 10).Value = 2000
$ws.Cells.Item(32, 11).Value = 5925
$ws.Cells.Item(32, 12).Value = 2000
$ws.Cells.Item(32, 13).Value = -5638
$ws.Cells.Item(32, 14).Value = -2574

# ARM row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(63, 8).Value = 2966.6667
$ws.Cells.Item(63, 9).Value = 2966.6667
$ws.Cells.Item(63, 11).Value = 2966.6667
$ws.Cells.Item(63, 13).Value = -2280.6667

# ARM row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(66, 8).Value = 2966.6667
$ws.Cells.Item(66, 9).Value = 2966.6667
$ws.Cells.Item(66, 11).Value = 14833.3335
$ws.Cells.Item(66, 13).Value = -11401.3335

# BSM row 26
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(26, 8).Value = 30471
$ws.Cells.Item(26, 9).Value = 30471
$ws.Cells.Item(26, 11).Value = 30471
$ws.Cells.Item(26, 13).Value = -30179

# BSM row 76
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(76, 8).Value = 0
$ws.Cells.Item(76, 10).Value = 0
$ws.Cells.Item(76, 12).Value = $null
$ws.Cells.Item(76, 14).Value = 0

# BSM row 79
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(79, 8).Value = 0
$ws.Cells.Item(79, 10).Value = 0
$ws.Cells.Item(79, 12).Value = $null
$ws.Cells.Item(79, 14).Value = 0

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 38501.75
$ws.Cells.Item(86, 9).Value = 2000
$ws.Cells.Item(86, 10).Value = 75003.5
$ws.Cells.Item(86, 11).Value = 2000
$ws.Cells.Item(86, 12).Value = 75003.5
$ws.Cells.Item(86, 13).Value = -877
$ws.Cells.Item(86, 14).Value = -77249.5

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(89, 8).Value = 38501.75
$ws.Cells.Item(89, 9).Value = 2000
$ws.Cells.Item(89, 10).Value = 75003.5
$ws.Cells.Item(89, 11).Value = 10000
$ws.Cells.Item(89, 12).Value = 375017.5
$ws.Cells.Item(89, 13).Value = -4384
$ws.Cells.Item(89, 14).Value = -386249.5

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2793.1738
$ws.Cells.Item(31, 9).Value = 1404.8572
$ws.Cells.Item(31, 11).Value = 1404.8572
$ws.Cells.Item(31, 13).Value = -1109.8572

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 2793.1738
$ws.Cells.Item(34, 9).Value = 1404.8572
$ws.Cells.Item(34, 11).Value = 1404.8572
$ws.Cells.Item(34, 13).Value = -1202.8572

# CRP row 45
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(45, 8).Value = 0
$ws.Cells.Item(45, 9).Value = 0
$ws.Cells.Item(45, 11).Value = 0
$ws.Cells.Item(45, 13).Value = $null

# CRP row 74
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(74, 8).Value = 10000
$ws.Cells.Item(74, 10).Value = 0
$ws.Cells.Item(74, 12).Value = 0
$ws.Cells.Item(74, 14).Value = $null

# CRP row 77
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(77, 8).Value = 10000
$ws.Cells.Item(77, 10).Value = 0
$ws.Cells.Item(77, 12).Value = 0
$ws.Cells.Item(77, 14).Value = $null

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 2500
$ws.Cells.Item(99, 9).Value = 0
$ws.Cells.Item(99, 10).Value = 2500
$ws.Cells.Item(99, 11).Value = 0
$ws.Cells.Item(99, 12).Value = $null
$ws.Cells.Item(99, 13).Value = 2500
$ws.Cells.Item(99, 14).Value = -5496

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(126, 8).Value = 2500
$ws.Cells.Item(126, 9).Value = 0
$ws.Cells.Item(126, 10).Value = 2500
$ws.Cells.Item(126, 11).Value = 0
$ws.Cells.Item(126, 12).Value = $null
$ws.Cells.Item(126, 13).Value = 7500
$ws.Cells.Item(126, 14).Value = -12440

# CUL row 25
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(25, 8).Value = 457.2

# CUL row 30
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(30, 8).Value = 457.2

# CUL row 39
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(39, 8).Value = 3939.4
$ws.Cells.Item(39, 10).Value = 4355.5454
$ws.Cells.Item(39, 12).Value = 13066.6362
$ws.Cells.Item(39, 14).Value = -13654.6362

# CUL row 55
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(55, 8).Value = 2679.4285
$ws.Cells.Item(55, 10).Value = 2834.3333
$ws.Cells.Item(55, 12).Value = 8502.999899999999
$ws.Cells.Item(55, 14).Value = -8856.999899999999

# CUL row 57
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(57, 8).Value = 100
$ws.Cells.Item(57, 9).Value = 100
$ws.Cells.Item(57, 11).Value = 300
$ws.Cells.Item(57, 13).Value = 259

# CUL row 59
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(59, 8).Value = 549.3333
$ws.Cells.Item(59, 9).Value = 749
$ws.Cells.Item(59, 10).Value = 150
$ws.Cells.Item(59, 11).Value = 2247
$ws.Cells.Item(59, 12).Value = 450
$ws.Cells.Item(59, 13).Value = -1707
$ws.Cells.Item(59, 14).Value = -1530

# CUL row 64
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(64, 8).Value = 1860
$ws.Cells.Item(64, 9).Value = 1825
$ws.Cells.Item(64, 10).Value = 2000
$ws.Cells.Item(64, 11).Value = 5475
$ws.Cells.Item(64, 12).Value = 6000
$ws.Cells.Item(64, 13).Value = -5205
$ws.Cells.Item(64, 14).Value = -6540

# CUL row 67
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(67, 8).Value = 1860
$ws.Cells.Item(67, 9).Value = 1825
$ws.Cells.Item(67, 10).Value = 2000
$ws.Cells.Item(67, 11).Value = 5475
$ws.Cells.Item(67, 12).Value = 6000
$ws.Cells.Item(67, 13).Value = -4539
$ws.Cells.Item(67, 14).Value = -7872

# CUL row 69
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(69, 8).Value = 24358.666
$ws.Cells.Item(69, 9).Value = 28537.5
$ws.Cells.Item(69, 11).Value = 85612.5
$ws.Cells.Item(69, 13).Value = -84801.5

# CUL row 72
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(72, 8).Value = 24358.666
$ws.Cells.Item(72, 9).Value = 28537.5
$ws.Cells.Item(72, 11).Value = 256837.5
$ws.Cells.Item(72, 13).Value = -252781.5

# CUL row 76
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(76, 8).Value = 850
$ws.Cells.Item(76, 9).Value = 850
$ws.Cells.Item(76, 11).Value = 2550
$ws.Cells.Item(76, 13).Value = -2167

# CUL row 79
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(79, 8).Value = 850
$ws.Cells.Item(79, 9).Value = 850
$ws.Cells.Item(79, 11).Value = 2550
$ws.Cells.Item(79, 13).Value = -1224

# CUL row 86
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(86, 8).Value = 19001.5

# CUL row 89
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(89, 8).Value = 19001.5

# CUL row 129
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(129, 8).Value = 3191.3333
$ws.Cells.Item(129, 10).Value = 3366.2222
$ws.Cells.Item(129, 12).Value = 10098.6666
$ws.Cells.Item(129, 14).Value = -20098.6666

# GSM row 31
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(31, 8).Value = 566.6667
$ws.Cells.Item(31, 9).Value = 566.6667
$ws.Cells.Item(31, 11).Value = 566.6667
$ws.Cells.Item(31, 13).Value = -274.6667

# GSM row 37
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(37, 8).Value = 566.6667
$ws.Cells.Item(37, 9).Value = 566.6667
$ws.Cells.Item(37, 11).Value = 566.6667
$ws.Cells.Item(37, 13).Value = -289.6667

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 2250
$ws.Cells.Item(126, 10).Value = 2500
$ws.Cells.Item(126, 12).Value = 7500
$ws.Cells.Item(126, 14).Value = -12440

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 2261.2666
$ws.Cells.Item(132, 9).Value = 1784.9166
$ws.Cells.Item(132, 10).Value = 4166.6665
$ws.Cells.Item(132, 11).Value = 5354.7498
$ws.Cells.Item(132, 12).Value = 12499.9995
$ws.Cells.Item(132, 13).Value = -2824.7498
$ws.Cells.Item(132, 14).Value = -17559.9995

# GSM row 136
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(136, 8).Value = 47499.5
$ws.Cells.Item(136, 10).Value = 47499.5
$ws.Cells.Item(136, 12).Value = 142498.5
$ws.Cells.Item(136, 14).Value = -147598.5

# GSM row 140
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(140, 8).Value = 115000
$ws.Cells.Item(140, 10).Value = 115000
$ws.Cells.Item(140, 12).Value = 115000
$ws.Cells.Item(140, 14).Value = -125360

# GSM row 141
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(141, 8).Value = 0
$ws.Cells.Item(141, 10).Value = 0
$ws.Cells.Item(141, 12).Value = $null
$ws.Cells.Item(141, 14).Value = 0

# LTW row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(82, 9).Value = 2237.5
$ws.Cells.Item(82, 10).Value = 1997.25
$ws.Cells.Item(82, 11).Value = 2237.5
$ws.Cells.Item(82, 12).Value = 1997.25
$ws.Cells.Item(82, 13).Value = -1876.5
$ws.Cells.Item(82, 14).Value = -2719.25

# LTW row 85
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(85, 9).Value = 2237.5
$ws.Cells.Item(85, 10).Value = 1997.25
$ws.Cells.Item(85, 11).Value = 2237.5
$ws.Cells.Item(85, 12).Value = 1997.25
$ws.Cells.Item(85, 13).Value = -989.5
$ws.Cells.Item(85, 14).Value = -4493.25

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 4769.8
$ws.Cells.Item(122, 9).Value = 4769.8
$ws.Cells.Item(122, 11).Value = 14309.4
$ws.Cells.Item(122, 13).Value = -11859.4

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 3433.889
$ws.Cells.Item(126, 9).Value = 2700.8572
$ws.Cells.Item(126, 10).Value = 5999.5
$ws.Cells.Item(126, 11).Value = 8102.571599999999
$ws.Cells.Item(126, 12).Value = 17998.5
$ws.Cells.Item(126, 13).Value = -5632.571599999999
$ws.Cells.Item(126, 14).Value = -22938.5

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 3213
$ws.Cells.Item(132, 9).Value = 3155.6
$ws.Cells.Item(132, 10).Value = 3500
$ws.Cells.Item(132, 11).Value = 9466.799999999999
$ws.Cells.Item(132, 12).Value = 10500
$ws.Cells.Item(132, 13).Value = -6936.799999999999
$ws.Cells.Item(132, 14).Value = -15560

# WVR row 140
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(140, 8).Value = 59997
$ws.Cells.Item(140, 10).Value = 59997
$ws.Cells.Item(140, 12).Value = 59997
$ws.Cells.Item(140, 14).Value = -70357
